# "block 2 and project plan"
# Fill in a planning header row, add a "Wer"/"Timeline" mini-table next to the
# requirement list, record who's doing what for item (5), move item (6)'s
# notes away, add a submission-deadline date, and re-highlight two items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 1: submission date note -----------------------------------
$ws.Range("A1").Value = "Abgabedatum 28.05."

# --- New header cells on row 2 (Wer / Timeline columns) -----------------
$ws.Range("C2").Value = "Wer"
$ws.Range("D2").Value = "Timeline"

# --- Item (1) "Web Scraping" note gets expanded + assignees -------------
$ws.Range("B5").Value = "Web Scraping von Börsennachrichten"
$ws.Range("C5").Value = "Michael (?), David (Cash), Sebastian (Yahoo Fianance)"

# --- Item (2) loses its old "Befehle suchen" / "DB iterieren" notes -----
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
# Item (2) is no longer highlighted in gold (theme7); matches item (3)/(4)
# style (theme9) instead.
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null

# --- Item (5) gets an assignee note --------------------------------------
$ws.Range("B9").Value = "siehe Unterricht Cirillo"
$ws.Range("C9").Value = "Michael"

# --- Item (6) gets an assignee note --------------------------------------
$ws.Range("B10").Value = "siehe Unterricht Spindler"

# --- Item (7) gets a timeline / due date ---------------------------------
$ws.Range("D11").NumberFormat = "d-mmm"
$ws.Range("D11").Value = (Get-Date -Year 2023 -Month 5 -Day 3 -Hour 0 -Minute 0 -Second 0).Date

# --- Item (9) is now highlighted in gold (theme7), like items (1)/(9 old) -
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null

# --- Block 2, item (2) gets highlighted green (theme9) -------------------
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null

# --- Block 2, item (4) gets an assignee note -----------------------------
$ws.Range("C19").Value = "David"

# --- Block 2, item (6) gets an assignee note -----------------------------
$ws.Range("C21").Value = "Sebastian"

# --- Resize the new "Wer"/"Timeline" columns to fit their content --------
$ws.Columns("B:C").AutoFit() | Out-Null

$excel.CutCopyMode = 0
Write-Output "done"
